# Updates the "cryptos" price/volume table (columns D and E, rows 2-51)
# with refreshed scrape values, per the Sat Jul 29 06:33:57 UTC 2023
# GitHub Actions run. Column D holds Price as text (it can contain
# multiple "." grouping separators, e.g. "29.373.64", so it is never a
# real number) and column E holds the 1h volume change percentage text
# (padded with spaces, e.g. "  +0.59%  "). For the handful of Price
# values that happen to look like a plain decimal number (e.g.
# "241.88"), the cell is briefly switched to Text format so Excel
# stores the digits verbatim instead of coercing them into a Number
# (which would also silently drop meaningful trailing zeros such as
# "0.07770" -> 0.0777); the style is then reset back to Normal so the
# cell's formatting matches the rest of the untouched column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.373.64'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.874.75'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.88'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07770'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.53%  '
$ws.Range('E10').Value = '  +1.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08444'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('D12').Value = '1.872.24'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.256'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7133'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').Value = '29.374.95'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.090'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008249'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.79'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').Value = '2.124.62'
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.777'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1594'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.063'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.513'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.425'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('E31').Value = '  +2.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.286'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('E33').Value = '  +3.48%  '
$ws.Range('E34').Value = '  +1.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.179'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7422'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.10%  '
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01871'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('D39').Value = '1.227.29'
$ws.Range('E39').Value = '  +4.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.731'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.522'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '110.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8892'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '2.021.92'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.814'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5216'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000123'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.450'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4321'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.37%  '
